# Added Filter for Customer in the Product Dependency Feature
#
# This script updates the "AddCustomer" sheet by adding an uppercase
# "filter" helper column (column D) that mirrors the customer names
# already present in column E (and updates the D1 header to the
# uppercase form of the E1 header). It also bumps the "Web Data" counter
# shown in column C of every sheet, and restores the selection/active
# sheet state left by the author after making the edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Update the shared "Web Data" counter (column C, row 1) - this
#    text is shared across all four sheets.
# ---------------------------------------------------------------
$orderPeriod = $wb.Worksheets.Item("OrderPeriod")
$addProductCategory = $wb.Worksheets.Item("AddProductCategory1")
$addCustomer = $wb.Worksheets.Item("AddCustomer")
$createOrder = $wb.Worksheets.Item("CreateOrder")

$orderPeriod.Range("C1").Value = "Web Data 49"
$addProductCategory.Range("C1").Value = "Web Data 49"
$addCustomer.Range("C1").Value = "Web Data 49"
$createOrder.Range("C1").Value = "Web Data 49"

# ---------------------------------------------------------------
# 2. AddCustomer sheet: add the new "filter" column D with the
#    uppercase customer names, mirroring column E.
# ---------------------------------------------------------------
$addCustomer.Range("D1").Value = "LEO HOLDER"
$addCustomer.Range("D2").Value = "WILLIAM"
$addCustomer.Range("D3").Value = "OLIVIA"
$addCustomer.Range("D4").Value = "ISABELLA"
$addCustomer.Range("D5").Value = "SOPHIA"
$addCustomer.Range("D6").Value = "JACOB"
$addCustomer.Range("D7").Value = "JAYDEN"

# Row 7 (Jayden) wraps text in column E, same as the rest of the
# header row - mirror that on the new D7 cell.
$addCustomer.Range("D7").WrapText = $true

# Match the width of the new column D to column E so the two
# columns line up as a single filter block.
$addCustomer.Columns.Item(4).ColumnWidth = 19.846666666666668

# Restore the author's final selection on the AddCustomer sheet.
$addCustomer.Range("E18").Select() | Out-Null

# ---------------------------------------------------------------
# 3. OrderPeriod sheet: selection moved back to C1 and it is no
#    longer the active/selected tab.
# ---------------------------------------------------------------
$orderPeriod.Range("C1").Select() | Out-Null

# ---------------------------------------------------------------
# 4. CreateOrder becomes the active sheet/tab, with the view
#    scrolled back to the top-left corner (A1).
# ---------------------------------------------------------------
$createOrder.Activate() | Out-Null
$createOrder.Range("A1").Select() | Out-Null
$createOrder.Range("I1").Select() | Out-Null
